$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell reference -> new displayed text, scraped from the refreshed
# coinranking.com snapshot (GitHub Actions cron job).
$updates = @(
    ,@("D2", "52.129.86")
    ,@("E2", "  +0.76%  ")
    ,@("D3", "2.897.62")
    ,@("E3", "  +3.48%  ")
    ,@("E4", "  +0.16%  ")
    ,@("D5", "351.63")
    ,@("E5", "  -0.44%  ")
    ,@("D6", "113.06")
    ,@("E6", "  +1.14%  ")
    ,@("E7", "  -0.24%  ")
    ,@("E8", "  +0.05%  ")
    ,@("D9", "0.621")
    ,@("E9", "  -0.61%  ")
    ,@("D10", "39.75")
    ,@("E10", "  -1.16%  ")
    ,@("E11", "  +0.55%  ")
    ,@("D12", "0.0858")
    ,@("E12", "  +2.49%  ")
    ,@("D13", "19.75")
    ,@("E13", "  -0.89%  ")
    ,@("D14", "7.68")
    ,@("E14", "  -0.99%  ")
    ,@("D15", "3.356.86")
    ,@("E15", "  +4.02%  ")
    ,@("D16", "2.913.79")
    ,@("E16", "  +4.91%  ")
    ,@("D17", "0.984")
    ,@("E17", "  +4.12%  ")
    ,@("D18", "52.217.60")
    ,@("E18", "  +0.99%  ")
    ,@("D19", "3.30")
    ,@("E19", "  +3.27%  ")
    ,@("D20", "7.57")
    ,@("E20", "  -0.34%  ")
    ,@("D21", "13.97")
    ,@("E21", "  +3.04%  ")
    ,@("D22", "0.0₃0971")
    ,@("E22", "  -0.04%  ")
    ,@("D23", "70.74")
    ,@("E23", "  +0.62%  ")
    ,@("D24", "268.37")
    ,@("E24", "  +0.41%  ")
    ,@("E25", "  +0.77%  ")
    ,@("D26", "0.178")
    ,@("E26", "  +11.72%  ")
    ,@("D27", "26.66")
    ,@("E27", "  +2.06%  ")
    ,@("D28", "1.00")
    ,@("E28", "  -0.03%  ")
    ,@("D29", "10.60")
    ,@("E29", "  +2.24%  ")
    ,@("D30", "0.102")
    ,@("E30", "  +14.98%  ")
    ,@("D31", "6.61")
    ,@("E31", "  +7.89%  ")
    ,@("D32", "37.39")
    ,@("E32", "  -4.01%  ")
    ,@("B33", "Toncoin")
    ,@("C33", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton")
    ,@("D33", "2.25")
    ,@("E33", "  -0.81%  ")
    ,@("B34", "RenderToken")
    ,@("C34", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr")
    ,@("D34", "6.20")
    ,@("E34", "  +11.83%  ")
    ,@("B35", "OKB")
    ,@("C35", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb")
    ,@("D35", "52.93")
    ,@("E35", "  +0.97%  ")
    ,@("D36", "0.0449")
    ,@("E36", "  -0.79%  ")
    ,@("E37", "  +0.06%  ")
    ,@("D38", "3.29")
    ,@("E38", "  +3.90%  ")
    ,@("D39", "18.79")
    ,@("E39", "  -0.34%  ")
    ,@("E40", "  +1.08%  ")
    ,@("D41", "2.72")
    ,@("E41", "  +8.20%  ")
    ,@("E42", "  +1.27%  ")
    ,@("D43", "23.05")
    ,@("E43", "  +5.45%  ")
    ,@("E44", "  +4.72%  ")
    ,@("D45", "119.45")
    ,@("E45", "  -0.70%  ")
    ,@("E46", "  -1.75%  ")
    ,@("D47", "2.167.96")
    ,@("E47", "  +3.07%  ")
    ,@("E48", "  +1.51%  ")
    ,@("D49", "0.261")
    ,@("E49", "  +19.58%  ")
    ,@("D50", "0.0343")
    ,@("E50", "  +10.69%  ")
    ,@("D51", "0.944")
    ,@("E51", "  -1.17%  ")
)

foreach ($pair in $updates) {
    $ref = $pair[0]
    $text = $pair[1]
    $cell = $ws.Range($ref)
    # Every D/E entry in this sheet is stored as literal text (already-
    # formatted strings, some of which look like plain numbers, e.g.
    # "3.30", "1.00"). Force text mode while writing so Excel does not
    # silently coerce them into numeric cells, then restore the default
    # style so no stray number formatting is left behind.
    if ($text -match "^[+-]?[0-9]+(\.[0-9]+)?$") {
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $text
    }
}
